$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell E1 - matches style of D1 (s="1"), numeric value 3
$ws.Range("E1").Value = 3
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats

# Data column E (rows 2-12), no special style, numeric values
$ws.Range("E2").Value = 0.002364
$ws.Range("E3").Value = 0.001306
$ws.Range("E4").Value = 0.000738
$ws.Range("E5").Value = 0.000559
$ws.Range("E6").Value = 0.000627
$ws.Range("E7").Value = 0.000816
$ws.Range("E8").Value = 0.000845
$ws.Range("E9").Value = 0.000963
$ws.Range("E10").Value = 0.001044
$ws.Range("E11").Value = 0.000394
$ws.Range("E12").Value = 0.000626
